$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptos price/volume table (rows 2-51) with the latest scrape
# values. Numeric-looking price strings are written with a leading
# apostrophe so Excel keeps them as text (matching the original inlineStr
# cells) instead of silently coercing them into floating point numbers.
$ws.Range('D2').Value = '62.016.52'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '3.420.68'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '''410.21'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').Value = '''129.99'
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('D7').Value = '''0.634'
$ws.Range('E7').Value = '  +6.97%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '''0.746'
$ws.Range('E9').Value = '  +10.25%  '
$ws.Range('D10').Value = '''0.145'
$ws.Range('E10').Value = '  +18.18%  '
$ws.Range('D11').Value = '''42.86'
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').Value = '''0.0000218'
$ws.Range('E12').Value = '  +66.55%  '
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').Value = '3.973.16'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('E15').Value = '  +6.37%  '
$ws.Range('D16').Value = '''21.07'
$ws.Range('E16').Value = '  +5.43%  '
$ws.Range('D17').Value = '3.429.32'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '''12.30'
$ws.Range('E18').Value = '  +11.42%  '
$ws.Range('E19').Value = '  +4.49%  '
$ws.Range('D20').Value = '62.003.38'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').Value = '''403.10'
$ws.Range('E21').Value = '  +27.51%  '
$ws.Range('D22').Value = '''89.91'
$ws.Range('E22').Value = '  +5.75%  '
$ws.Range('D23').Value = '''3.19'
$ws.Range('E23').Value = '  -1.00%  '
$ws.Range('D24').Value = '''13.21'
$ws.Range('E24').Value = '  +3.03%  '
$ws.Range('D25').Value = '''3.26'
$ws.Range('E25').Value = '  +4.07%  '
$ws.Range('D26').Value = '''32.69'
$ws.Range('E26').Value = '  +9.76%  '
$ws.Range('D27').Value = '''8.65'
$ws.Range('E27').Value = '  +5.08%  '
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').Value = '''7.61'
$ws.Range('E29').Value = '  -1.24%  '
$ws.Range('D30').Value = '''0.120'
$ws.Range('E30').Value = '  +3.49%  '
$ws.Range('D31').Value = '''2.70'
$ws.Range('E31').Value = '  -1.28%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '''0.172'
$ws.Range('E32').Value = '  -1.56%  '
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').Value = '''11.92'
$ws.Range('E33').Value = '  +4.24%  '
$ws.Range('D34').Value = '''43.40'
$ws.Range('E34').Value = '  +0.88%  '
$ws.Range('E35').Value = '  +0.68%  '
$ws.Range('E36').Value = '  +2.31%  '
$ws.Range('D37').Value = '''54.42'
$ws.Range('E37').Value = '  +4.43%  '
$ws.Range('D38').Value = '''0.998'
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('D39').Value = '''3.39'
$ws.Range('E39').Value = '  -1.55%  '
$ws.Range('E40').Value = '  +7.42%  '
$ws.Range('E41').Value = '  -2.53%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = '''143.07'
$ws.Range('E42').Value = '  +3.53%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').Value = '''0.312'
$ws.Range('E43').Value = '  +6.86%  '
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('D45').Value = '''4.09'
$ws.Range('E45').Value = '  +2.20%  '
$ws.Range('D46').Value = '''2.42'
$ws.Range('E46').Value = '  +9.05%  '
$ws.Range('D47').Value = '''16.69'
$ws.Range('E47').Value = '  -1.08%  '
$ws.Range('D48').Value = '''21.83'
$ws.Range('E48').Value = '  +1.49%  '
$ws.Range('D49').Value = '2.129.00'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').Value = '''2.37'
$ws.Range('E50').Value = '  +3.45%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.132'
$ws.Range('E51').Value = '  +17.16%  '
